$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.180.47"
$ws.Range("E2").Value = "  +7.60%  "
$ws.Range("D3").Value = "2.582.93"
$ws.Range("E3").Value = "  +9.89%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'504.76"
$ws.Range("E5").Value = "  +7.00%  "
$ws.Range("D6").Value = "'156.53"
$ws.Range("E6").Value = "  +8.71%  "
$ws.Range("D7").Value = "'0.637"
$ws.Range("E7").Value = "  +26.85%  "
$ws.Range("D8").Value = "'0.997"
$ws.Range("E8").Value = "  -0.29%  "
$ws.Range("D9").Value = "2.580.01"
$ws.Range("E9").Value = "  +9.75%  "
$ws.Range("D10").Value = "'6.19"
$ws.Range("E10").Value = "  +14.74%  "
$ws.Range("E11").Value = "  +6.86%  "
$ws.Range("E12").Value = "  +6.71%  "
$ws.Range("E13").Value = "  +1.90%  "
$ws.Range("D14").Value = "2.995.70"
$ws.Range("E14").Value = "  +8.62%  "
$ws.Range("D15").Value = "59.105.87"
$ws.Range("E15").Value = "  +7.35%  "
$ws.Range("D16").Value = "'21.72"
$ws.Range("E16").Value = "  +8.53%  "
$ws.Range("E17").Value = "  +5.39%  "
$ws.Range("D18").Value = "2.571.86"
$ws.Range("E18").Value = "  +9.37%  "
$ws.Range("E19").Value = "  +5.24%  "
$ws.Range("D20").Value = "'336.25"
$ws.Range("D21").Value = "'10.34"
$ws.Range("E21").Value = "  +8.15%  "
$ws.Range("D22").Value = "'6.05"
$ws.Range("E22").Value = "  +8.28%  "
$ws.Range("E23").Value = "  +0.81%  "
$ws.Range("D24").Value = "'59.81"
$ws.Range("E24").Value = "  +7.19%  "
$ws.Range("E25").Value = "  +6.11%  "
$ws.Range("E26").Value = "  +8.46%  "
$ws.Range("D27").Value = "2.677.64"
$ws.Range("E27").Value = "  +9.24%  "
$ws.Range("D28").Value = "'0.997"
$ws.Range("E28").Value = "  -0.35%  "
$ws.Range("D29").Value = "0.0₃0826"
$ws.Range("E29").Value = "  +9.91%  "
$ws.Range("D30").Value = "'7.39"
$ws.Range("E30").Value = "  +3.22%  "
$ws.Range("D31").Value = "'0.997"
$ws.Range("E31").Value = "  -0.32%  "
$ws.Range("D32").Value = "'157.43"
$ws.Range("E32").Value = "  +6.24%  "
$ws.Range("D33").Value = "'19.29"
$ws.Range("E33").Value = "  +7.54%  "
$ws.Range("E34").Value = "  +7.21%  "
$ws.Range("D35").Value = "'5.49"
$ws.Range("E35").Value = "  +9.63%  "
$ws.Range("E36").Value = "  +10.54%  "
$ws.Range("E37").Value = "  +9.62%  "
$ws.Range("D38").Value = "'0.850"
$ws.Range("E38").Value = "  +3.57%  "
$ws.Range("E39").Value = "  +12.10%  "
$ws.Range("E40").Value = "  +8.67%  "
$ws.Range("D41").Value = "'35.16"
$ws.Range("E41").Value = "  +5.02%  "
$ws.Range("D42").Value = "'292.48"
$ws.Range("E42").Value = "  +14.99%  "
$ws.Range("D43").Value = "'0.103"
$ws.Range("E43").Value = "  +8.68%  "
$ws.Range("D44").Value = "'0.624"
$ws.Range("E44").Value = "  +8.56%  "
$ws.Range("D45").Value = "'0.0562"
$ws.Range("E45").Value = "  +7.18%  "
$ws.Range("D46").Value = "'0.995"
$ws.Range("E46").Value = "  -0.28%  "
$ws.Range("D47").Value = "'0.755"
$ws.Range("E47").Value = "  +20.31%  "
$ws.Range("D48").Value = "'19.29"
$ws.Range("E48").Value = "  +15.29%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'4.83"
$ws.Range("E49").Value = "  +9.36%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0235"
$ws.Range("E50").Value = "  +6.84%  "
$ws.Range("E51").Value = "  +0.88%  "
